$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge the two runs "...endangered species" + "." into a single run.
#    A Find/Replace across the combined text collapses the matched range
#    into one run (taking on the first run's formatting), exactly as the
#    diff shows.
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "API calls were made to retrieve NY Times articles of the endangered species.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "API calls were made to retrieve NY Times articles of the endangered species.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2. Remove the _GoBack bookmark from its old location (the empty
#    paragraph right after "Transform"). It is a hidden bookmark so it
#    doesn't show up in Bookmarks.Count/iteration, but it can still be
#    reached directly by name.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3. Add the two new paragraphs right after the "Temperature: ..." bullet:
#      - a blank paragraph
#      - a paragraph with the new "App has been created..." text, which
#        also carries the (re-created) _GoBack bookmark
# ---------------------------------------------------------------------

# Locate the "Temperature: ..." paragraph and the existing blank
# paragraph that immediately follows it.
$tempPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Temperature: Country, Average Temperature*") {
        $tempPara = $p
        break
    }
}

$existingBlank = $tempPara.Next()

# Insert the new "text" paragraph before the pre-existing blank line -
# this keeps its paragraph formatting clean (no numbering list carried
# over from the "Temperature:" bullet).
$existingBlank.Range.InsertParagraphBefore() | Out-Null
$textPara = $existingBlank.Previous()
$textPara.Range.Text = "App has been created to visually show the data in html. Link for downloading the temperature data by extant country for each queried animal is available."

# Add the _GoBack bookmark to the very end of that paragraph (collapsed,
# matching its original empty-range placement).
$bmRange = $textPara.Range.Duplicate
$bmRange.MoveEnd(1, -1) | Out-Null
$bmRange.Collapse(0) | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Insert the blank paragraph before the text paragraph. Type then delete
# a placeholder character so the engine doesn't leave a stray empty run
# behind in the serialized XML.
$textPara.Range.InsertParagraphBefore() | Out-Null
$blankPara = $textPara.Previous()
$blankPara.Range.Text = "x"
$blankRange = $blankPara.Range.Duplicate
$blankRange.MoveEnd(1, -1) | Out-Null
$blankRange.Delete() | Out-Null
